# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5,6).Value = 175
$ws.Cells.Item(6,6).Value = 827
$ws.Cells.Item(7,6).Value = 4229
$ws.Cells.Item(8,6).Value = 4229
$ws.Cells.Item(12,6).Value = 6151
$ws.Cells.Item(13,6).Value = 67
$ws.Cells.Item(14,6).Value = 467
$ws.Cells.Item(14,7).Value = "已售罄"
$ws.Cells.Item(17,6).Value = 169
$ws.Cells.Item(19,6).Value = 9269
$ws.Cells.Item(21,6).Value = 2498
$ws.Cells.Item(22,6).Value = 196
$ws.Cells.Item(24,6).Value = 2470
$ws.Cells.Item(27,6).Value = 1978
$ws.Cells.Item(30,6).Value = 333
$ws.Cells.Item(32,6).Value = 45
$ws.Cells.Item(35,6).Value = 73
$ws.Cells.Item(37,6).Value = 1224
$ws.Cells.Item(38,6).Value = 1223
$ws.Cells.Item(40,6).Value = 101
$ws.Cells.Item(42,6).Value = 1558
$ws.Cells.Item(43,6).Value = 2562
$ws.Cells.Item(45,6).Value = 932
$ws.Cells.Item(46,6).Value = 309
$ws.Cells.Item(47,6).Value = 1255

# --- 演出 (sheet2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5,7).Value = "不可售"
$ws.Cells.Item(9,6).Value = 12
$ws.Cells.Item(22,6).Value = 85

# --- 本地生活 (sheet3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4,6).Value = 104

# --- 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(10,2).NumberFormat = "@"
$ws.Cells.Item(10,2).Value = "2024-08-10"
$ws.Cells.Item(10,3).Value = "北京· 人气声优 樱川惠 专场活动"
$ws.Cells.Item(10,4).Value = "北京展览馆 北京展览馆"
$ws.Cells.Item(10,5).Value = "2024.08.10 13:10-08.10 16:30"
$ws.Cells.Item(10,6).Value = 175
$ws.Cells.Item(10,7).Value = 468
$ws.Cells.Item(10,8).Value = "https://show.bilibili.com/platform/detail.html?id=87897"
$ws.Cells.Item(10,9).Value = "//i1.hdslb.com/bfs/openplatform/202406/5m9c4FJT1719210386733.png"
$ws.Cells.Item(11,3).Value = "北京·AINI二次元派对【免票展会】"
$ws.Cells.Item(11,4).Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws.Cells.Item(11,5).Value = "2024.08.10 10:00-08.10 16:00"
$ws.Cells.Item(11,6).Value = 827
$ws.Cells.Item(11,7).Value = 50
$ws.Cells.Item(11,8).Value = "https://show.bilibili.com/platform/detail.html?id=89601"
$ws.Cells.Item(11,9).Value = "//i1.hdslb.com/bfs/openplatform/202407/eIryW6Up1721208870214.jpeg"
$ws.Cells.Item(12,3).Value = "北京·GOJO超次元动漫游戏嘉年华15th"
$ws.Cells.Item(12,4).Value = "小关路39号 北投购物公园"
$ws.Cells.Item(12,5).Value = "2024.08.10 09:20-08.11 17:00"
$ws.Cells.Item(12,6).Value = 4229
$ws.Cells.Item(12,7).Value = 6.6
$ws.Cells.Item(12,8).Value = "https://show.bilibili.com/platform/detail.html?id=85223"
$ws.Cells.Item(12,9).Value = "//i0.hdslb.com/bfs/openplatform/202407/DlsfmegR1720613188484.jpeg"
$ws.Cells.Item(13,3).Value = "北京·【七夕专场-告白多巴胺】“我心永恒”唯美经典电影作品音乐会"
$ws.Cells.Item(13,4).Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$ws.Cells.Item(13,5).Value = "2024.08.10 19:30-08.10 21:00"
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 90
$ws.Cells.Item(13,8).Value = "https://show.bilibili.com/platform/detail.html?id=89478"
$ws.Cells.Item(13,9).Value = "//i2.hdslb.com/bfs/openplatform/202407/zVFK9v7b1720088899389.png"
$ws.Cells.Item(17,3).Value = "北京·我心永恒——唯美英文经典歌曲七夕演唱会"
$ws.Cells.Item(17,4).Value = "复兴门内大街49号 民族宫大剧院"
$ws.Cells.Item(17,5).Value = "2024.08.10 19:30-08.10 21:30"
$ws.Cells.Item(17,6).Value = 12
$ws.Cells.Item(17,7).Value = 99
$ws.Cells.Item(17,8).Value = "https://show.bilibili.com/platform/detail.html?id=87228"
$ws.Cells.Item(17,9).Value = "//i2.hdslb.com/bfs/openplatform/202406/kFhqQvyX1718162118549.jpeg"
$ws.Cells.Item(18,3).Value = "北京·梦次元动漫展M30"
$ws.Cells.Item(18,4).Value = "北京展览馆 北京展览馆"
$ws.Cells.Item(18,5).Value = "2024.08.10 10:00-08.11 17:00"
$ws.Cells.Item(18,6).Value = 6151
$ws.Cells.Item(18,7).Value = 80
$ws.Cells.Item(18,8).Value = "https://show.bilibili.com/platform/detail.html?id=83828"
$ws.Cells.Item(18,9).Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"
$ws.Cells.Item(19,3).Value = "北京·狐妖小红娘专题聚会【免票活动】"
$ws.Cells.Item(19,4).Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws.Cells.Item(19,5).Value = "2024.08.10 14:00-08.10 18:00"
$ws.Cells.Item(19,6).Value = 67
$ws.Cells.Item(19,7).Value = 58
$ws.Cells.Item(19,8).Value = "https://show.bilibili.com/platform/detail.html?id=90238"
$ws.Cells.Item(19,9).Value = "//i1.hdslb.com/bfs/openplatform/202408/mL8ytYCG1722578125040.jpeg"
$ws.Cells.Item(20,6).Value = 169
$ws.Cells.Item(22,6).Value = 9269
$ws.Cells.Item(24,6).Value = 2498
$ws.Cells.Item(25,6).Value = 196
$ws.Cells.Item(27,6).Value = 2470
$ws.Cells.Item(30,6).Value = 1978
$ws.Cells.Item(33,6).Value = 333
$ws.Cells.Item(36,6).Value = 73
$ws.Cells.Item(38,6).Value = 1223
$ws.Cells.Item(40,6).Value = 101
$ws.Cells.Item(42,6).Value = 1558
$ws.Cells.Item(43,6).Value = 2562
$ws.Cells.Item(44,6).Value = 932
$ws.Cells.Item(45,6).Value = 309
$ws.Cells.Item(48,6).Value = 1255
$ws.Cells.Item(50,6).Value = 85
$ws.Cells.Item(51,6).Value = 85
